# Assignment5.docx — "week 5 and 8 update"
#
# 1) Part 1 paragraph: the runs holding "lightly inaccurate" and ", "
#    are merged into a single run "lightly inaccurate, ".
# 2) The "94.3%" run gets a yellow highlight, and the "_GoBack" bookmark
#    (Word's "last edit" marker) moves to sit right after it — which
#    also removes it from its old spot after the "large difference"
#    sentence near the end of the document, since a document can only
#    have one bookmark with a given name.

$d = $word.ActiveDocument

# --- 1) merge "lightly inaccurate" + ", " into one run ------------------
$rng = $d.Content
$rng.Find.Execute(
    "lightly inaccurate, ",   # FindText
    $true,                    # MatchCase
    $false,                   # MatchWholeWord
    $false,                   # MatchWildcards
    $false,                   # MatchSoundsLike
    $false,                   # MatchAllWordForms
    $true,                    # Forward
    1,                        # Wrap (wdFindContinue)
    $false,                   # Format
    "lightly inaccurate, ",   # ReplaceWith
    2                         # Replace (wdReplaceAll)
) | Out-Null

# --- 2) highlight "94.3%" and move the _GoBack bookmark after it -------
$rng = $d.Content
$rng.Find.Execute("94.3%", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rng.Font.HighlightColorIndex = 7   # wdYellow
$bmRange = $d.Range($rng.End, $rng.End)
$d.Bookmarks.Add("_GoBack", $bmRange) | Out-Null
